# Apply the "Wrote API overview" protocol entry and update the saved
# selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 68: new protocol entry dated 2013-04-16 (Excel serial date 41380)
# describing the work done ("Wrote API overview"). A68 already carries the
# date number format (style index 4) in the template, so only the values
# need to be written; Excel will append the new text to the shared string
# table automatically.
$ws.Cells.Item(68, 1).Value = 41380
$ws.Cells.Item(68, 2).Value = "Wrote API overview"

# Move the sheet's saved selection from B67 to B64.
$ws.Range("B64").Select()
